$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Staff data update: add honorifics to two newly added staff names
# (shared strings for the surrounding unique_id / photo-path cells are left untouched)
$ws.Range("A18").Value = "Mr. SHANMUGAPRIYAN S"
$ws.Range("A19").Value = "Mrs. NISHA DEVI N"

# Selection / scroll position ends up at A20 (just past the last data row)
# with the view scrolled back so column A is visible again.
$ws.Range("A20").Select()
